$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("site")
$ws.Activate()

# Insert a new column before column B (shifts Watershed..Other right by one)
$ws.Columns("B:B").Insert()

# New column header for the inserted "Reach Number" field
$ws.Range("B2").Value = "Reach Number"

# Fix the casing of the site-name column header to match the rest of the table
$ws.Range("A2").Value = "Site"

# Leave the selection where the author left it after adding the column
$ws.Range("C9").Select()
